# Update "人气"/heat-count values (column F) on the "展览" and "全部类型" sheets
# to match the newly generated data output (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet ---
$wsExhibition.Range("F3").Value  = 2712
$wsExhibition.Range("F7").Value  = 2296
$wsExhibition.Range("F9").Value  = 218
$wsExhibition.Range("F12").Value = 555
$wsExhibition.Range("F18").Value = 9235
$wsExhibition.Range("F20").Value = 7178
$wsExhibition.Range("F21").Value = 11728
$wsExhibition.Range("F27").Value = 2596
$wsExhibition.Range("F30").Value = 2549
$wsExhibition.Range("F34").Value = 916
$wsExhibition.Range("F37").Value = 531

# --- 全部类型 (All types) sheet ---
$wsAllTypes.Range("F5").Value  = 2712
$wsAllTypes.Range("F9").Value  = 2296
$wsAllTypes.Range("F13").Value = 218
$wsAllTypes.Range("F16").Value = 555
$wsAllTypes.Range("F22").Value = 9235
$wsAllTypes.Range("F24").Value = 7178
$wsAllTypes.Range("F25").Value = 11728
$wsAllTypes.Range("F33").Value = 2596
$wsAllTypes.Range("F46").Value = 531

Write-Output "Updated 20 cells across 展览 and 全部类型 sheets"
